$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 81
$ws.Range("F6").Value = 2626
$ws.Range("F7").Value = 60
$ws.Range("F9").Value = 528
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 1529
$ws.Range("F12").Value = 3
$ws.Range("F14").Value = 629
$ws.Range("F15").Value = 1507
$ws.Range("F16").Value = 1353
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 534
$ws.Range("F19").Value = 3793
$ws.Range("F21").Value = 3297
$ws.Range("F22").Value = 756
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 2186
$ws.Range("F26").Value = 303
$ws.Range("F28").Value = 21
$ws.Range("F29").Value = 1160
$ws.Range("F32").Value = 1037
$ws.Range("F33").Value = 1029

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 9
$ws.Range("F12").Value = 100
$ws.Range("F17").Value = 115
$ws.Range("F18").Value = 247
$ws.Range("F20").Value = 475

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 528
$ws.Range("F6").Value = 103

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 81
$ws.Range("F9").Value = 528
$ws.Range("F10").Value = 103
$ws.Range("F11").Value = 2626
$ws.Range("F12").Value = 2626
$ws.Range("F14").Value = 60
$ws.Range("F19").Value = 528
$ws.Range("F22").Value = 1529
$ws.Range("F23").Value = 9
$ws.Range("F25").Value = 1507
$ws.Range("F26").Value = 100
$ws.Range("F27").Value = 1353
$ws.Range("F28").Value = 17
$ws.Range("F29").Value = 535
$ws.Range("F31").Value = 3793
$ws.Range("F33").Value = 3297
$ws.Range("F34").Value = 756
$ws.Range("F35").Value = 2186
$ws.Range("F37").Value = 303
$ws.Range("F39").Value = 21
$ws.Range("F40").Value = 1160
$ws.Range("F42").Value = 115
$ws.Range("F43").Value = 247
$ws.Range("F45").Value = 475
$ws.Range("F48").Value = 1037
$ws.Range("F49").Value = 1029
